$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the bottom-border formatting from the current last data row (29) onto
# row 25, which will become the new last row of the table once rows 26:29
# are removed below.
$ws.Range("B29:J29").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 21-25 get refreshed with a new "part 1" batch: the same five workers
# from the 2507 period, now carrying the new 2508 period/mora values.
$ws.Range("C21").Value = "73169606"
$ws.Range("D21").Value = "CESAR EDUARDO MATOS FRANCO"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("C22").Value = "78382447"
$ws.Range("D22").Value = "DAVID MANUEL FERIA PEREZ"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("C23").Value = "1047502589"
$ws.Range("D23").Value = "INGRIT BANESA FLOREZ CORREA"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500

$ws.Range("C24").Value = "1143402667"
$ws.Range("D24").Value = "MANUEL CORREA GARCIA"
$ws.Range("E24").Value = "2508"
$ws.Range("F24").Value = 56940
$ws.Range("G24").Value = 1423500

$ws.Range("C25").Value = "1007976046"
$ws.Range("D25").Value = "JAIDER RAFAEL CORREA GARCIA"
$ws.Range("E25").Value = "2508"
$ws.Range("F25").Value = 56940
$ws.Range("G25").Value = 1423500

# The old periods 2506/2505/2504/2503/2502 rows (26-29) are no longer part
# of the account statement; remove them and let everything below shift up.
$ws.Rows("26:29").Delete()

# Refresh the summary counters: 5 workers, 2 periods (2507 & 2508), and the
# new total overdue amount.
$ws.Range("E11").Value = 569400
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 2

# Column D ("Nombre Trabajador") auto-fits to its content; the longest name
# left in the table is shorter than before, so the column narrows.
$ws.Columns("D:D").AutoFit()
